{"js": "// Minor changes to modelling report:\n// 1) \"... A good computer can maybe do a 6x6 board. It is interesting ...\"\n//    -> \"... A good computer can maybe do a 5x5 board in less than an hour. It is interesting ...\"\n// 2) \"... 2 rooks can put a king in checkmate on a 4x4 board, as well as a 5x5 board. ...\"\n//    -> \"... 2 rooks can put a king in stalemate on a 4x4 board, as well as a 5x5 board. ...\"\n\nconst body = context.document.body;\n\n// Change 1: board size/timing update.\nconst boardSizeResults = body.search(\"6x6 board.\", { matchCase: true, matchWholeWord: false });\nboardSizeResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < boardSizeResults.items.length; i++) {\n  boardSizeResults.items[i].insertText(\n    \"5x5 board in less than an hour.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// Change 2: checkmate -> stalemate for the 2-rooks example.\nconst rookExampleResults = body.search(\"in checkmate on a 4x4\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nrookExampleResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < rookExampleResults.items.length; i++) {\n  rookExampleResults.items[i].insertText(\n    \"in stalemate on a 4x4\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Minor changes to modelling report:\n# 1) \"... A good computer can maybe do a 6x6 board. It is interesting ...\"\n#    -> \"... A good computer can maybe do a 5x5 board in less than an hour. It is interesting ...\"\n# 2) \"... 2 rooks can put a king in checkmate on a 4x4 board, as well as a 5x5 board. ...\"\n#    -> \"... 2 rooks can put a king in stalemate on a 4x4 board, as well as a 5x5 board. ...\"\n\n$d = $word.ActiveDocument\n\n# Change 1: board size/timing update.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"6x6 board.\"\n$find1.Replacement.Text = \"5x5 board in less than an hour.\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# Change 2: checkmate -> stalemate for the 2-rooks example.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"in checkmate on a 4x4\"\n$find2.Replacement.Text = \"in stalemate on a 4x4\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
